$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1099.8868
$ws.Range("I15").Value = 1099.8868
$ws.Range("K15").Value = 3299.6604
$ws.Range("M15").Value = -3130.6604
$ws.Range("H19").Value = 1145.7858
$ws.Range("I19").Value = 541
$ws.Range("J19").Value = 1750.5714
$ws.Range("K19").Value = 541
$ws.Range("L19").Value = 1750.5714
$ws.Range("M19").Value = -366
$ws.Range("N19").Value = -2100.5714
$ws.Range("H28").Value = 894.36365
$ws.Range("I28").Value = 894.36365
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 894.36365
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -409.36365
$ws.Range("N28").ClearContents()
$ws.Range("H76").Value = 2530544
$ws.Range("I76").Value = 3273011
$ws.Range("J76").Value = 6156.4
$ws.Range("K76").Value = 3273011
$ws.Range("L76").Value = 6156.4
$ws.Range("M76").Value = -3272696
$ws.Range("N76").Value = -6786.4
$ws.Range("H79").Value = 2530544
$ws.Range("I79").Value = 3273011
$ws.Range("J79").Value = 6156.4
$ws.Range("K79").Value = 3273011
$ws.Range("L79").Value = 6156.4
$ws.Range("M79").Value = -3271919
$ws.Range("N79").Value = -8340.4
$ws.Range("H80").Value = 8226.666999999999
$ws.Range("I80").Value = 2466.6667
$ws.Range("J80").Value = 9666.666999999999
$ws.Range("K80").Value = 7400.000100000001
$ws.Range("L80").Value = 29000.001
$ws.Range("M80").Value = -6402.000100000001
$ws.Range("N80").Value = -30996.001
$ws.Range("H83").Value = 8226.666999999999
$ws.Range("I83").Value = 2466.6667
$ws.Range("J83").Value = 9666.666999999999
$ws.Range("K83").Value = 22200.0003
$ws.Range("L83").Value = 87000.003
$ws.Range("M83").Value = -17208.0003
$ws.Range("N83").Value = -96984.003
$ws.Range("H88").Value = 3351.5715
$ws.Range("I88").Value = 1999.5
$ws.Range("K88").Value = 1999.5
$ws.Range("M88").Value = -1593.5
$ws.Range("H91").Value = 3351.5715
$ws.Range("I91").Value = 1999.5
$ws.Range("K91").Value = 1999.5
$ws.Range("M91").Value = -595.5
$ws.Range("H113").Value = 7820
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 7820
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 7820
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -14328
$ws.Range("H116").Value = 5633.3687
$ws.Range("J116").Value = 5684.067
$ws.Range("L116").Value = 5684.067
$ws.Range("N116").Value = -12568.067
$ws.Range("H132").Value = 3339.8333
$ws.Range("I132").Value = 3140.3655
$ws.Range("J132").Value = 4636.375
$ws.Range("K132").Value = 9421.0965
$ws.Range("L132").Value = 13909.125
$ws.Range("M132").Value = -6891.0965
$ws.Range("N132").Value = -18969.125
$ws.Range("H135").Value = 1408.0834
$ws.Range("I135").Value = 1194.9524
$ws.Range("K135").Value = 10754.5716
$ws.Range("M135").Value = -8219.571599999999
$ws.Range("H137").Value = 88604.766
$ws.Range("I137").Value = 114534.625
$ws.Range("J137").Value = 5629.2
$ws.Range("K137").Value = 343603.875
$ws.Range("L137").Value = 16887.6
$ws.Range("M137").Value = -341053.875
$ws.Range("N137").Value = -21987.6
$ws.Range("H138").Value = 3538.3015
$ws.Range("J138").Value = 4665.4185
$ws.Range("L138").Value = 13996.2555
$ws.Range("N138").Value = -24276.2555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 63791.438
$ws.Range("I2").Value = 72632.14
$ws.Range("J2").Value = 1906.5
$ws.Range("K2").Value = 72632.14
$ws.Range("L2").Value = 1906.5
$ws.Range("M2").Value = -72519.14
$ws.Range("N2").Value = -2132.5
$ws.Range("H4").Value = 196.33333
$ws.Range("I4").Value = 196.33333
$ws.Range("K4").Value = 196.33333
$ws.Range("M4").Value = -80.33332999999999
$ws.Range("H32").Value = 4382.6665
$ws.Range("I32").Value = 2491.603
$ws.Range("J32").Value = 11150.685
$ws.Range("K32").Value = 2491.603
$ws.Range("L32").Value = 11150.685
$ws.Range("M32").Value = -2204.603
$ws.Range("N32").Value = -11724.685
$ws.Range("H74").Value = 86071.47
$ws.Range("I74").Value = 40597.824
$ws.Range("J74").Value = 202281.89
$ws.Range("K74").Value = 40597.824
$ws.Range("L74").Value = 202281.89
$ws.Range("M74").Value = -39723.824
$ws.Range("N74").Value = -204029.89
$ws.Range("H77").Value = 86071.47
$ws.Range("I77").Value = 40597.824
$ws.Range("J77").Value = 202281.89
$ws.Range("K77").Value = 202989.12
$ws.Range("L77").Value = 1011409.45
$ws.Range("M77").Value = -198621.12
$ws.Range("N77").Value = -1020145.45
$ws.Range("H108").Value = 37342
$ws.Range("J108").Value = 37342
$ws.Range("L108").Value = 37342
$ws.Range("N108").Value = -45022
$ws.Range("H110").Value = 1292.4546
$ws.Range("I110").Value = 1191.7
$ws.Range("K110").Value = 1191.7
$ws.Range("M110").Value = 853.3
$ws.Range("H116").Value = 63791.438
$ws.Range("I116").Value = 72632.14
$ws.Range("J116").Value = 1906.5
$ws.Range("K116").Value = 72632.14
$ws.Range("L116").Value = 1906.5
$ws.Range("M116").Value = -70338.14
$ws.Range("N116").Value = -6494.5
$ws.Range("H122").Value = 3082.5
$ws.Range("I122").Value = 3049.5
$ws.Range("K122").Value = 9148.5
$ws.Range("M122").Value = -6698.5
$ws.Range("H132").Value = 2428.6082
$ws.Range("I132").Value = 1660.6666
$ws.Range("J132").Value = 3846.3462
$ws.Range("K132").Value = 4981.9998
$ws.Range("L132").Value = 11539.0386
$ws.Range("M132").Value = -2451.9998
$ws.Range("N132").Value = -16599.0386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 63791.438
$ws.Range("I3").Value = 72632.14
$ws.Range("J3").Value = 1906.5
$ws.Range("K3").Value = 72632.14
$ws.Range("L3").Value = 1906.5
$ws.Range("M3").Value = -72518.14
$ws.Range("N3").Value = -2134.5
$ws.Range("H68").Value = 74518
$ws.Range("J68").Value = 74518
$ws.Range("L68").Value = 74518
$ws.Range("N68").Value = -76140
$ws.Range("H71").Value = 74518
$ws.Range("J71").Value = 74518
$ws.Range("L71").Value = 223554
$ws.Range("N71").Value = -231666
$ws.Range("H86").Value = 6347.8887
$ws.Range("I86").Value = 6940.636
$ws.Range("K86").Value = 6940.636
$ws.Range("M86").Value = -5817.636
$ws.Range("H89").Value = 6347.8887
$ws.Range("I89").Value = 6940.636
$ws.Range("K89").Value = 34703.18
$ws.Range("M89").Value = -29087.18
$ws.Range("H94").Value = 6746.0303
$ws.Range("I94").Value = 973.8929000000001
$ws.Range("J94").Value = 39070
$ws.Range("K94").Value = 973.8929000000001
$ws.Range("L94").Value = 39070
$ws.Range("M94").Value = -522.8929000000001
$ws.Range("N94").Value = -39972
$ws.Range("H107").Value = 2598.75
$ws.Range("I107").Value = 2201.4546
$ws.Range("K107").Value = 2201.4546
$ws.Range("M107").Value = -281.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 240.82353
$ws.Range("I7").Value = 233.06667
$ws.Range("J7").Value = 299
$ws.Range("K7").Value = 233.06667
$ws.Range("L7").Value = 299
$ws.Range("M7").Value = -120.06667
$ws.Range("N7").Value = -525
$ws.Range("H16").Value = 2396.4285
$ws.Range("I16").Value = 1639.4445
$ws.Range("K16").Value = 1639.4445
$ws.Range("M16").Value = -1352.4445
$ws.Range("H31").Value = 20327.02
$ws.Range("I31").Value = 1922.8462
$ws.Range("K31").Value = 1922.8462
$ws.Range("M31").Value = -1627.8462
$ws.Range("H34").Value = 20327.02
$ws.Range("I34").Value = 1922.8462
$ws.Range("K34").Value = 1922.8462
$ws.Range("M34").Value = -1720.8462
$ws.Range("H59").Value = 41442.23
$ws.Range("J59").Value = 40729.082
$ws.Range("L59").Value = 40729.082
$ws.Range("N59").Value = -43019.082
$ws.Range("H74").Value = 49999
$ws.Range("J74").Value = 49999
$ws.Range("L74").Value = 49999
$ws.Range("N74").Value = -51747
$ws.Range("H77").Value = 49999
$ws.Range("J77").Value = 49999
$ws.Range("L77").Value = 149997
$ws.Range("N77").Value = -158733
$ws.Range("H99").Value = 3732.4666
$ws.Range("I99").Value = 3274
$ws.Range("J99").Value = 4256.4287
$ws.Range("K99").Value = 3274
$ws.Range("L99").Value = 4256.4287
$ws.Range("M99").Value = -1776
$ws.Range("N99").Value = -7252.4287
$ws.Range("H107").Value = 27028386
$ws.Range("I107").Value = 1349.3572
$ws.Range("K107").Value = 1349.3572
$ws.Range("M107").Value = 570.6428000000001
$ws.Range("H113").Value = 2396.4285
$ws.Range("I113").Value = 1639.4445
$ws.Range("K113").Value = 1639.4445
$ws.Range("M113").Value = 530.5554999999999
$ws.Range("H126").Value = 3732.4666
$ws.Range("I126").Value = 3274
$ws.Range("J126").Value = 4256.4287
$ws.Range("K126").Value = 9822
$ws.Range("L126").Value = 12769.2861
$ws.Range("M126").Value = -7352
$ws.Range("N126").Value = -17709.2861
$ws.Range("H132").Value = 59298.625
$ws.Range("I132").Value = 3619
$ws.Range("J132").Value = 84607.55
$ws.Range("K132").Value = 10857
$ws.Range("L132").Value = 253822.65
$ws.Range("M132").Value = -8327
$ws.Range("N132").Value = -258882.65

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 337.9697
$ws.Range("I2").Value = 123.454544
$ws.Range("K2").Value = 740.727264
$ws.Range("M2").Value = -627.727264
$ws.Range("H11").Value = 74725.10000000001
$ws.Range("J11").Value = 132998.25
$ws.Range("L11").Value = 398994.75
$ws.Range("N11").Value = -399274.75
$ws.Range("H12").Value = 40466.91
$ws.Range("I12").Value = 148234.33
$ws.Range("K12").Value = 444702.99
$ws.Range("M12").Value = -444529.99
$ws.Range("H34").Value = 627.7778
$ws.Range("I34").Value = 537.5
$ws.Range("K34").Value = 1612.5
$ws.Range("M34").Value = -1528.5
$ws.Range("H39").Value = 5443
$ws.Range("I39").Value = 1329
$ws.Range("K39").Value = 3987
$ws.Range("M39").Value = -3693
$ws.Range("H55").Value = 86355.836
$ws.Range("I55").Value = 895.7143
$ws.Range("J55").Value = 206000
$ws.Range("K55").Value = 2687.1429
$ws.Range("L55").Value = 618000
$ws.Range("M55").Value = -2510.1429
$ws.Range("N55").Value = -618354
$ws.Range("H60").Value = 733.3333
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 733.3333
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 2199.9999
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -2701.9999
$ws.Range("H68").Value = 1228.125
$ws.Range("J68").Value = 3047
$ws.Range("L68").Value = 9141
$ws.Range("N68").Value = -10763
$ws.Range("H70").Value = 736.6667
$ws.Range("I70").Value = 736.6667
$ws.Range("K70").Value = 2210.0001
$ws.Range("M70").Value = -1895.0001
$ws.Range("H71").Value = 1228.125
$ws.Range("J71").Value = 3047
$ws.Range("L71").Value = 27423
$ws.Range("N71").Value = -35535
$ws.Range("H73").Value = 736.6667
$ws.Range("I73").Value = 736.6667
$ws.Range("K73").Value = 2210.0001
$ws.Range("M73").Value = -1118.0001
$ws.Range("H92").Value = 1468.1111
$ws.Range("J92").Value = 2247
$ws.Range("L92").Value = 6741
$ws.Range("N92").Value = -9237
$ws.Range("H122").Value = 1005.5
$ws.Range("I122").Value = 946.75
$ws.Range("J122").Value = 1093.625
$ws.Range("K122").Value = 8520.75
$ws.Range("L122").Value = 9842.625
$ws.Range("M122").Value = -6070.75
$ws.Range("N122").Value = -14742.625
$ws.Range("H137").Value = 1858.5555
$ws.Range("I137").Value = 1459.1538
$ws.Range("K137").Value = 4377.4614
$ws.Range("M137").Value = 722.5385999999999
$ws.Range("H139").Value = 22728346
$ws.Range("I139").Value = 27778392
$ws.Range("K139").Value = 83335176
$ws.Range("M139").Value = -83330036
$ws.Range("H140").Value = 2062.2144
$ws.Range("I140").Value = 1759.3846
$ws.Range("K140").Value = 5278.1538
$ws.Range("M140").Value = -98.15380000000005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 9398.5
$ws.Range("J54").Value = 9398.5
$ws.Range("L54").Value = 9398.5
$ws.Range("N54").Value = -10178.5
$ws.Range("H80").Value = 22820260
$ws.Range("I80").Value = 45456810
$ws.Range("K80").Value = 45456810
$ws.Range("M80").Value = -45455812
$ws.Range("H83").Value = 22820260
$ws.Range("I83").Value = 45456810
$ws.Range("K83").Value = 227284050
$ws.Range("M83").Value = -227279058
$ws.Range("H113").Value = 3563.375
$ws.Range("I113").Value = 2929.5715
$ws.Range("K113").Value = 2929.5715
$ws.Range("M113").Value = -759.5715
$ws.Range("H122").Value = 209769.23
$ws.Range("I122").Value = 250538.33
$ws.Range("J122").Value = 50399.184
$ws.Range("K122").Value = 751614.99
$ws.Range("L122").Value = 151197.552
$ws.Range("M122").Value = -749164.99
$ws.Range("N122").Value = -156097.552
$ws.Range("H132").Value = 2682.7144
$ws.Range("I132").Value = 2502.186
$ws.Range("J132").Value = 3976.5
$ws.Range("K132").Value = 7506.558000000001
$ws.Range("L132").Value = 11929.5
$ws.Range("M132").Value = -4976.558000000001
$ws.Range("N132").Value = -16989.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 5748.5
$ws.Range("I30").Value = 5748.5
$ws.Range("K30").Value = 5748.5
$ws.Range("M30").Value = -5640.5
$ws.Range("H40").Value = 6951.4
$ws.Range("I40").Value = 5713.2354
$ws.Range("K40").Value = 5713.2354
$ws.Range("M40").Value = -5577.2354
$ws.Range("H46").Value = 4511.25
$ws.Range("J46").Value = 5504
$ws.Range("L46").Value = 5504
$ws.Range("N46").Value = -5880
$ws.Range("H61").Value = 7304.9443
$ws.Range("I61").Value = 7518.8237
$ws.Range("K61").Value = 7518.8237
$ws.Range("M61").Value = -7316.8237
$ws.Range("H113").Value = 7304.9443
$ws.Range("I113").Value = 7518.8237
$ws.Range("K113").Value = 7518.8237
$ws.Range("M113").Value = -5348.8237
$ws.Range("H132").Value = 8395.447
$ws.Range("I132").Value = 8857.360000000001
$ws.Range("J132").Value = 7507.154
$ws.Range("K132").Value = 26572.08
$ws.Range("L132").Value = 22521.462
$ws.Range("M132").Value = -24042.08
$ws.Range("N132").Value = -27581.462
$ws.Range("H133").Value = 106211.86
$ws.Range("J133").Value = 106211.86
$ws.Range("L133").Value = 106211.86
$ws.Range("N133").Value = -111271.86

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8816.117
$ws.Range("I62").Value = 8320.625
$ws.Range("K62").Value = 8320.625
$ws.Range("M62").Value = -7696.625
$ws.Range("H65").Value = 8816.117
$ws.Range("I65").Value = 8320.625
$ws.Range("K65").Value = 41603.125
$ws.Range("M65").Value = -38483.125
$ws.Range("H96").Value = 5559.2
$ws.Range("I96").Value = 4300
$ws.Range("K96").Value = 4300
$ws.Range("M96").Value = -2927
$ws.Range("H100").Value = 2220
$ws.Range("I100").Value = 1604.4445
$ws.Range("K100").Value = 3208.889
$ws.Range("M100").Value = -2667.889
$ws.Range("H107").Value = 1500
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4500
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2580
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 1378.7693
$ws.Range("I113").Value = 847.75
$ws.Range("K113").Value = 2543.25
$ws.Range("M113").Value = -373.25
$ws.Range("H132").Value = 19932.908
$ws.Range("I132").Value = 3887.0977
$ws.Range("J132").Value = 70538.92
$ws.Range("K132").Value = 11661.2931
$ws.Range("L132").Value = 211616.76
$ws.Range("M132").Value = -9131.293099999999
$ws.Range("N132").Value = -216676.76
